# Replace the label "Valor médio unitário na silvicultura" with
# "Preço médio recebido na silvicultura" in column B of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = $ws.Range("B1:B734")
$colB.Replace("Valor médio unitário na silvicultura", "Preço médio recebido na silvicultura", 1, 1, $false, $false, $false)
